$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 614-615; all data currently at row 614 onward
# (through the old last row 713) shifts down by two rows, to 616-715.
$ws.Rows("614:615").Insert()

# Populate the two newly-inserted rows with the new "Región Metropolitana"
# entries (same constant columns as the surrounding block: A, B, C, E, F,
# G, H, N, Q, R).
$ws.Range("A614").Value2 = 3
$ws.Range("B614").Value = "Femacal de La Calera"
$ws.Range("C614").Value = "Coquimbo"
$ws.Range("D614").Value2 = 44995
$ws.Range("E614").Value2 = 5
$ws.Range("F614").Value2 = 100112028
$ws.Range("G614").Value = "Sandia"
$ws.Range("H614").Value = "Sin especificar"
$ws.Range("I614").Value = "Primera"
$ws.Range("J614").Value2 = 380
$ws.Range("K614").Value2 = 3000
$ws.Range("L614").Value2 = 3000
$ws.Range("M614").Value2 = 3000
$ws.Range("N614").Value = "$/unidad"
$ws.Range("O614").Value = "Región Metropolitana"
$ws.Range("P614").Value2 = 3000
$ws.Range("Q614").Value2 = 1
$ws.Range("R614").Value = "Hortaliza"

$ws.Range("A615").Value2 = 3
$ws.Range("B615").Value = "Femacal de La Calera"
$ws.Range("C615").Value = "Coquimbo"
$ws.Range("D615").Value2 = 44995
$ws.Range("E615").Value2 = 5
$ws.Range("F615").Value2 = 100112028
$ws.Range("G615").Value = "Sandia"
$ws.Range("H615").Value = "Sin especificar"
$ws.Range("I615").Value = "Segunda"
$ws.Range("J615").Value2 = 250
$ws.Range("K615").Value2 = 2000
$ws.Range("L615").Value2 = 2000
$ws.Range("M615").Value2 = 2000
$ws.Range("N615").Value = "$/unidad"
$ws.Range("O615").Value = "Región Metropolitana"
$ws.Range("P615").Value2 = 2000
$ws.Range("Q615").Value2 = 1
$ws.Range("R615").Value = "Hortaliza"
